# Fruta / hortaliza, semanal
# Insert a new weekly record at row 73 (pushing the existing rows 73-74
# down to 74-75) and populate it with the latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 73:74 down to 74:75, leaving a blank row 73 for the new entry.
$ws.Rows.Item(73).Insert()

$ws.Range("A73").Value = 8
$ws.Range("B73").Value = "Terminal La Palmera de La Serena"
$ws.Range("C73").Value = "Coquimbo"
$ws.Range("D73").Value = 44628
$ws.Range("E73").Value = 4
$ws.Range("F73").Value = 100112030
$ws.Range("G73").Value = "Poroto granado"
$ws.Range("H73").Value = "Sin especificar"
$ws.Range("I73").Value = "Primera"
$ws.Range("J73").Value = 560
$ws.Range("K73").Value = 29000
$ws.Range("L73").Value = 30000
$ws.Range("M73").Value = 29500
$ws.Range("N73").Value = "$/malla 25 kilos"
$ws.Range("O73").Value = "Provincia del Elquí"
$ws.Range("P73").Value = 1180
$ws.Range("Q73").Value = 25
$ws.Range("R73").Value = "Hortaliza"
